$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDT")

# Update SAMPLE1b (row 21) feedback text: was about the "second" chord / Nummer 2,
# now about the "third" chord / Nummer 3, and wording tweaked.
$ws.Range("B21").Value = "In diesem Beispiel wurde der dritte Akkord verändert. Die richtige Antwort wäre also **Nummer 3**. Es folgen nun zwei Übungsfragen."
$ws.Range("C21").Value = " Here, the third chord was different, so the correct answer would have been **number 3**. Now you will see two practice questions."

# Widen column B to fit the longer practice text (stored width ends up at 126.5).
$ws.Columns.Item(2).ColumnWidth = 125.66666666666667

# Restore the view/selection state to match what was saved.
$ws.Range("C21").Select()
$excel.ActiveWindow.ScrollRow = 11
